# Commit: "Commit atualizacao tesouro dia 26/06"
#
# The GRAFICO sheet keeps a rolling day-by-day table: row 3 holds the
# date for each AJ/AK column pair, row 4 holds the VENDA (sell) rate
# (Taxa de Rendimento) / price (Preco Unitario), row 5 holds the same
# for COMPRA (buy). The next open slot for the day's data was the
# AJ/AK pair (previously blank). Rows 6-8 already carry formulas that
# reference these cells ($B$6/AK4, AJ4-based interest accrual, etc.),
# so once the new values are written Excel's automatic recalculation
# fills in the dependent results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GRAFICO")
$ws.Activate()

# Row 3: date header for the new column -> 26/jun/2019 (serial 43642)
$ws.Range("AJ3").Value = 43642
$ws.Range("AK3").Value = 43642

# Row 4: VENDA (sell) - Taxa de Rendimento (% a.a.) / Preco Unitario
$ws.Range("AJ4").Value = 3.85
$ws.Range("AK4").Value = 1775.24

# Row 5: COMPRA (buy) - Taxa de Rendimento (% a.a.) / Preco Unitario
$ws.Range("AJ5").Value = 3.73
$ws.Range("AK5").Value = 1808.03

# Recalculate so all dependent formulas (AJ6/AK6, AJ7/AK7, AV6/AV7/AV8,
# and the chart-feeding cells) pick up the new inputs.
$excel.CalculateFullRebuild()

# Reflect the cell the author ended up on after entering the day's data.
$ws.Range("AK4").Select()
